# Updates the "想去人数" (want-to-go count) figures in column F across the
# three worksheets that contain this column, per the scraped-data refresh
# commit "Update gh-pages to output generated at 456a3b4".
# Sheet order in the workbook: 1=展览 (Exhibition), 2=演出 (Performance),
# 3=本地生活 (Local Life, unaffected), 4=全部类型 (All Types).

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition)
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 19
$ws.Range("F3").Value = 28
$ws.Range("F5").Value = 5158
$ws.Range("F6").Value = 5158
$ws.Range("F7").Value = 119
$ws.Range("F8").Value = 133
$ws.Range("F9").Value = 515
$ws.Range("F11").Value = 1164
$ws.Range("F12").Value = 717
$ws.Range("F13").Value = 5028
$ws.Range("F15").Value = 64
$ws.Range("F16").Value = 80
$ws.Range("F17").Value = 226
$ws.Range("F18").Value = 233
$ws.Range("F20").Value = 246
$ws.Range("F21").Value = 3798
$ws.Range("F23").Value = 40
$ws.Range("F24").Value = 3696
$ws.Range("F25").Value = 178
$ws.Range("F26").Value = 172
$ws.Range("F28").Value = 219
$ws.Range("F29").Value = 237
$ws.Range("F30").Value = 205
$ws.Range("F31").Value = 105
$ws.Range("F32").Value = 110
$ws.Range("F36").Value = 6575
$ws.Range("F37").Value = 1053
$ws.Range("F38").Value = 494
$ws.Range("F39").Value = 98
$ws.Range("F40").Value = 970
$ws.Range("F41").Value = 58
$ws.Range("F42").Value = 1341
$ws.Range("F43").Value = 158
$ws.Range("F44").Value = 661
$ws.Range("F46").Value = 2251
$ws.Range("F49").Value = 771
$ws.Range("F50").Value = 913

# Sheet 2: 演出 (Performance)
$ws = $wb.Worksheets.Item(2)
$ws.Range("F8").Value = 50
$ws.Range("F9").Value = 85
$ws.Range("F14").Value = 4
$ws.Range("F24").Value = 807

# Sheet 4: 全部类型 (All Types)
$ws = $wb.Worksheets.Item(4)
$ws.Range("F4").Value = 19
$ws.Range("F6").Value = 28
$ws.Range("F8").Value = 5158
$ws.Range("F9").Value = 5158
$ws.Range("F10").Value = 119
$ws.Range("F11").Value = 50
$ws.Range("F12").Value = 133
$ws.Range("F13").Value = 85
$ws.Range("F14").Value = 515
$ws.Range("F15").Value = 1164
$ws.Range("F16").Value = 717
$ws.Range("F17").Value = 5028
$ws.Range("F19").Value = 64
$ws.Range("F20").Value = 80
$ws.Range("F21").Value = 226
$ws.Range("F22").Value = 233
$ws.Range("F24").Value = 246
$ws.Range("F25").Value = 3798
$ws.Range("F26").Value = 3696
$ws.Range("F27").Value = 178
$ws.Range("F28").Value = 172
$ws.Range("F29").Value = 219
$ws.Range("F30").Value = 237
$ws.Range("F31").Value = 205
$ws.Range("F32").Value = 105
$ws.Range("F33").Value = 110
$ws.Range("F37").Value = 6575
$ws.Range("F38").Value = 1053
$ws.Range("F39").Value = 494
$ws.Range("F41").Value = 98
$ws.Range("F42").Value = 970
$ws.Range("F43").Value = 1341
$ws.Range("F44").Value = 158
$ws.Range("F45").Value = 661
$ws.Range("F46").Value = 2251
$ws.Range("F48").Value = 771
$ws.Range("F49").Value = 913
